# Auto-generated Excel COM-interop script to apply scheduled price-data update
# to the Kujata_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H32").Value = 2150
$ws.Range("I32").Value = 966.6667
$ws.Range("J32").Value = 3333.3333
$ws.Range("K32").Value = 966.6667
$ws.Range("L32").Value = 3333.3333
$ws.Range("M32").Value = -640.6667
$ws.Range("N32").Value = -3985.3333

$ws.Range("H33").Value = 466.57693
$ws.Range("I33").Value = 450.8
$ws.Range("J33").Value = 519.1667
$ws.Range("K33").Value = 450.8
$ws.Range("L33").Value = 519.1667
$ws.Range("M33").Value = -221.8
$ws.Range("N33").Value = -977.1667

$ws.Range("H113").Value = 2927.5
$ws.Range("I113").Value = 2812
$ws.Range("K113").Value = 2812
$ws.Range("M113").Value = 442

$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 4635.6294
$ws.Range("I32").Value = 4429.365
$ws.Range("K32").Value = 4429.365
$ws.Range("M32").Value = -4142.365

$ws.Range("H45").Value = 1299.75
$ws.Range("I45").Value = 1234.4286
$ws.Range("K45").Value = 1234.4286
$ws.Range("M45").Value = -857.4286

$ws.Range("H61").Value = 37037916
$ws.Range("I61").Value = 41667480
$ws.Range("J61").Value = 1392.6666
$ws.Range("K61").Value = 41667480
$ws.Range("L61").Value = 1392.6666
$ws.Range("M61").Value = -41667268
$ws.Range("N61").Value = -1816.6666

$ws.Range("H131").Value = 50613
$ws.Range("J131").Value = 50613
$ws.Range("L131").Value = 50613
$ws.Range("N131").Value = -60693

$ws.Range("H132").Value = 2683.75
$ws.Range("I132").Value = 2160.45
$ws.Range("J132").Value = 3992
$ws.Range("K132").Value = 6481.349999999999
$ws.Range("L132").Value = 11976
$ws.Range("M132").Value = -3951.349999999999
$ws.Range("N132").Value = -17036

$ws.Range("H136").Value = 37037916
$ws.Range("I136").Value = 41667480
$ws.Range("J136").Value = 1392.6666
$ws.Range("K136").Value = 125002440
$ws.Range("L136").Value = 4177.9998
$ws.Range("M136").Value = -124999890
$ws.Range("N136").Value = -9277.9998


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H80").Value = 466.64706
$ws.Range("J80").Value = 547
$ws.Range("L80").Value = 547
$ws.Range("N80").Value = -2543

$ws.Range("H83").Value = 466.64706
$ws.Range("J83").Value = 547
$ws.Range("L83").Value = 2735
$ws.Range("N83").Value = -12719

$ws.Range("H86").Value = 2625.5715
$ws.Range("I86").Value = 2718.182
$ws.Range("J86").Value = 2523.7
$ws.Range("K86").Value = 2718.182
$ws.Range("L86").Value = 2523.7
$ws.Range("M86").Value = -1595.182
$ws.Range("N86").Value = -4769.7

$ws.Range("H89").Value = 2625.5715
$ws.Range("I89").Value = 2718.182
$ws.Range("J89").Value = 2523.7
$ws.Range("K89").Value = 13590.91
$ws.Range("L89").Value = 12618.5
$ws.Range("M89").Value = -7974.91
$ws.Range("N89").Value = -23850.5

$ws.Range("H105").Value = 58824532
$ws.Range("I105").Value = 58824532
$ws.Range("K105").Value = 58824532
$ws.Range("M105").Value = -58822785

$ws.Range("H107").Value = 1502.4546
$ws.Range("I107").Value = 1000.1429
$ws.Range("J107").Value = 2381.5
$ws.Range("K107").Value = 1000.1429
$ws.Range("L107").Value = 2381.5
$ws.Range("M107").Value = 919.8570999999999
$ws.Range("N107").Value = -6221.5

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 62501010
$ws.Range("I16").Value = 76923980
$ws.Range("J16").Value = 1450
$ws.Range("K16").Value = 76923980
$ws.Range("L16").Value = 1450
$ws.Range("M16").Value = -76923693
$ws.Range("N16").Value = -2024

$ws.Range("H22").Value = 70335.2
$ws.Range("J22").Value = 87870.125
$ws.Range("L22").Value = 87870.125
$ws.Range("N22").Value = -88570.125

$ws.Range("H99").Value = 1403.6364
$ws.Range("I99").Value = 1375
$ws.Range("J99").Value = 1438
$ws.Range("K99").Value = 1375
$ws.Range("L99").Value = 1438
$ws.Range("M99").Value = 123
$ws.Range("N99").Value = -4434

$ws.Range("H113").Value = 62501010
$ws.Range("I113").Value = 76923980
$ws.Range("J113").Value = 1450
$ws.Range("K113").Value = 76923980
$ws.Range("L113").Value = 1450
$ws.Range("M113").Value = -76921810
$ws.Range("N113").Value = -5790

$ws.Range("H126").Value = 1403.6364
$ws.Range("I126").Value = 1375
$ws.Range("J126").Value = 1438
$ws.Range("K126").Value = 4125
$ws.Range("L126").Value = 4314
$ws.Range("M126").Value = -1655
$ws.Range("N126").Value = -9254

$ws.Range("H140").Value = 37750
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 37750
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 37750
$ws.Range("N140").Value = -48110
$ws.Range("M140").ClearContents()


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H3").Value = 8824.040000000001
$ws.Range("I3").Value = 2954.6924
$ws.Range("J3").Value = 15182.5
$ws.Range("K3").Value = 8864.0772
$ws.Range("L3").Value = 45547.5
$ws.Range("M3").Value = -8752.0772
$ws.Range("N3").Value = -45771.5

$ws.Range("H56").Value = 5744.857
$ws.Range("I56").Value = 5744.857
$ws.Range("K56").Value = 5744.857
$ws.Range("M56").Value = -5214.857

$ws.Range("H97").Value = 927.55554
$ws.Range("I97").Value = 802
$ws.Range("J97").Value = 1084.5
$ws.Range("K97").Value = 2406
$ws.Range("L97").Value = 3253.5
$ws.Range("M97").Value = -1910
$ws.Range("N97").Value = -4245.5

$ws.Range("H131").Value = 37039060
$ws.Range("J131").Value = 2322.3044
$ws.Range("L131").Value = 6966.9132
$ws.Range("N131").Value = -17046.9132

$ws.Range("H139").Value = 1532.3611
$ws.Range("I139").Value = 1426.7273
$ws.Range("J139").Value = 1698.3572
$ws.Range("K139").Value = 4280.1819
$ws.Range("L139").Value = 5095.071599999999
$ws.Range("M139").Value = 859.8181000000004
$ws.Range("N139").Value = -15375.0716

$ws.Range("H140").Value = 23857.938
$ws.Range("J140").Value = 3398.5938
$ws.Range("L140").Value = 10195.7814
$ws.Range("N140").Value = -20555.7814


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1275.5
$ws.Range("I22").Value = 1400
$ws.Range("K22").Value = 1400
$ws.Range("M22").Value = -1105

$ws.Range("H27").Value = 1275.5
$ws.Range("I27").Value = 1400
$ws.Range("K27").Value = 1400
$ws.Range("M27").Value = -1293

$ws.Range("H46").Value = 8600
$ws.Range("J46").Value = 8600
$ws.Range("L46").Value = 8600
$ws.Range("N46").Value = -8976

$ws.Range("H61").Value = 1090
$ws.Range("I61").Value = 1141.8572
$ws.Range("J61").Value = 999.25
$ws.Range("K61").Value = 1141.8572
$ws.Range("L61").Value = 999.25
$ws.Range("M61").Value = -939.8571999999999
$ws.Range("N61").Value = -1403.25

$ws.Range("H68").Value = 1827.8572
$ws.Range("I68").Value = 1816.6666
$ws.Range("J68").Value = 1895
$ws.Range("K68").Value = 1816.6666
$ws.Range("L68").Value = 1895
$ws.Range("M68").Value = -1067.6666
$ws.Range("N68").Value = -3393

$ws.Range("H71").Value = 1827.8572
$ws.Range("I71").Value = 1816.6666
$ws.Range("J71").Value = 1895
$ws.Range("K71").Value = 9083.333000000001
$ws.Range("L71").Value = 9475
$ws.Range("M71").Value = -5339.333000000001
$ws.Range("N71").Value = -16963

$ws.Range("H113").Value = 1090
$ws.Range("I113").Value = 1141.8572
$ws.Range("J113").Value = 999.25
$ws.Range("K113").Value = 1141.8572
$ws.Range("L113").Value = 999.25
$ws.Range("M113").Value = 1028.1428
$ws.Range("N113").Value = -5339.25

$ws.Range("H122").Value = 17858698
$ws.Range("I122").Value = 25001226
$ws.Range("J122").Value = 2376.25
$ws.Range("K122").Value = 75003678
$ws.Range("L122").Value = 7128.75
$ws.Range("M122").Value = -75001228
$ws.Range("N122").Value = -12028.75


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H96").Value = 2485.5386
$ws.Range("I96").Value = 1789.7778
$ws.Range("K96").Value = 1789.7778
$ws.Range("M96").Value = -416.7778000000001

